# Generate Report for Handoff
# - Status flips from "Handed back: in sync with en-US" to "Ready for handoff"
# - Timestamps for the Overview / de-de "Handoff" sheets bump to 09:05:16
# - Timestamp for the zh-cn "Handoff" sheet bumps to 09:05:12
# - The now-shorter status text lets the Status column(s) narrow

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update status text "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Update Latest HO Xliff Generate Date / Latest Handoff Datetime timestamps
$overview.Range("G2").Value = "2016-08-20 09:05:16"
$dede.Range("H2").Value = "2016-08-20 09:05:16"
$zhcn.Range("H2").Value = "2016-08-20 09:05:12"

# Narrow the Status columns to fit the new, shorter text
$overview.Range("E:F").ColumnWidth = 16.3333333333333
$zhcn.Range("C:C").ColumnWidth = 16.3333333333333
$dede.Range("C:C").ColumnWidth = 16.3333333333333
